$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New signup row (row 4) matching the 15 header columns:
# timestamp, firstName, lastName, playedBefore, experienceLevel, playedClub,
# clubName, gender, hasDisability, location, email, phone, position, goal, whyJoin
$ws.Cells.Item(4, 1).Value  = "2025-07-08T14:41:16.912Z"
$ws.Cells.Item(4, 2).Value  = "jghnu67"
$ws.Cells.Item(4, 3).Value  = "hgju7"
$ws.Cells.Item(4, 4).Value  = $true
$ws.Cells.Item(4, 5).Value  = "Intermediate"
$ws.Cells.Item(4, 6).Value  = $true
$ws.Cells.Item(4, 7).Value  = "jghmuih"
$ws.Cells.Item(4, 8).Value  = "Male"
$ws.Cells.Item(4, 9).Value  = $false
$ws.Cells.Item(4, 10).Value = "jmhkui6trfdd"
$ws.Cells.Item(4, 11).Value = "hgngu@example.com"

# Phone number must be stored as TEXT (shared string), not a number -
# a plain .Value assignment of a purely-numeric string gets auto-coerced
# to a numeric cell, so build it as a text formula and paste back as a value.
$phoneCell = $ws.Cells.Item(4, 12)
$phoneCell.Formula = "=""12306978493"""
$phoneCell.Copy()
$phoneCell.PasteSpecial(-4163)

$ws.Cells.Item(4, 13).Value = "stricker"
$ws.Cells.Item(4, 14).Value = "7u7u7ijmhn"
$ws.Cells.Item(4, 15).Value = "sumn"
